# Updated cryptos list values per the authoritative diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Map of cell reference -> new text value. All target cells hold plain
# text (coin names, URLs, price strings, and padded percentage strings),
# so each cell's number format is forced to Text ('@') before assignment
# and then reset to the 'Normal' style so no stray formatting is left
# behind (this avoids Excel auto-converting numeric-looking strings such
# as "1.001" or "50.05" into real numbers).
$updates = [ordered]@{
    'D2' = '23.089.92'
    'E2' = '  -3.30%  '
    'D3' = '1.603.87'
    'E3' = '  -2.82%  '
    'E4' = '  -0.07%  '
    'D5' = '1.001'
    'E5' = '  -0.01%  '
    'D6' = '301.47'
    'E6' = '  -2.97%  '
    'D7' = '0.3781'
    'E7' = '  -2.96%  '
    'D8' = '0.3656'
    'E8' = '  -4.40%  '
    'D9' = '50.05'
    'E9' = '  -2.53%  '
    'D10' = '1.267'
    'E10' = '  -5.58%  '
    'B11' = 'BinanceUSD'
    'C11' = 'https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd'
    'D11' = '1.001'
    'E11' = '  +0.00%  '
    'B12' = 'Dogecoin'
    'C12' = 'https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge'
    'D12' = '0.08142'
    'E12' = '  -3.63%  '
    'D13' = '22.87'
    'E13' = '  -4.26%  '
    'D14' = '6.605'
    'E14' = '  -5.62%  '
    'E15' = '  -4.39%  '
    'D16' = '7.389'
    'E16' = '  -7.86%  '
    'D17' = '1.601.80'
    'E17' = '  -3.04%  '
    'D18' = '92.18'
    'E18' = '  -2.07%  '
    'D19' = '0.06887'
    'E19' = '  -1.33%  '
    'D20' = '18.29'
    'E20' = '  -6.31%  '
    'D21' = '6.578'
    'E21' = '  -5.33%  '
    'E22' = '  +0.03%  '
    'E23' = '  -3.82%  '
    'D24' = '23.084.47'
    'E24' = '  -3.34%  '
    'D25' = '2.364'
    'E25' = '  -3.28%  '
    'D26' = '2.797'
    'E26' = '  -4.46%  '
    'D27' = '21.14'
    'E27' = '  -3.95%  '
    'D28' = '150.50'
    'D29' = '5.271'
    'E29' = '  -2.35%  '
    'D30' = '134.80'
    'E30' = '  -1.81%  '
    'D31' = '2.364'
    'E31' = '  -4.74%  '
    'D32' = '6.866'
    'E32' = '  -11.29%  '
    'D33' = '1.778.48'
    'E33' = '  -2.79%  '
    'D34' = '0.9596'
    'E34' = '  -3.43%  '
    'D35' = '0.07693'
    'E35' = '  -5.48%  '
    'D36' = '10.50'
    'E36' = '  -1.71%  '
    'D37' = '6.304'
    'E37' = '  -5.30%  '
    'D38' = '0.02719'
    'E38' = '  -6.35%  '
    'D40' = '0.08914'
    'E40' = '  -2.10%  '
    'E41' = '  -3.23%  '
    'D42' = '0.7084'
    'E42' = '  -6.16%  '
    'D43' = '12.67'
    'E43' = '  -6.32%  '
    'D44' = '15.36'
    'E44' = '  -6.62%  '
    'D45' = '0.6646'
    'E45' = '  -4.06%  '
    'D46' = '2.339'
    'E46' = '  -4.19%  '
    'D47' = '1.000'
    'E47' = '  -0.03%  '
    'D48' = '4.005'
    'E48' = '  -2.56%  '
    'D49' = '132.47'
    'E49' = '  -0.60%  '
    'D50' = '1.248'
    'E50' = '  +1.87%  '
    'D51' = '0.07938'
    'E51' = '  -3.93%  '
}

foreach ($ref in $updates.Keys) {
    $cell = $ws.Range($ref)
    $cell.NumberFormat = '@'
    $cell.Value = $updates[$ref]
    $cell.Style = 'Normal'
}
